# PAM_STIM_2024 — slides were made full-screen, which shifted several of the
# PowerPoint slide references used to drive this story. Update the "Img"
# column (C) so the stimulus rows point at the renumbered slides, then leave
# the selection where the author left off (cell C19, just past the last
# used row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new slides were inserted right after "Diapositive6a", so the old
# "Diapositive6b"/"Diapositive6c" placeholders are renamed to the slides
# that now actually hold that content.
$ws.Range("C8").Value = "Diapositive7"
$ws.Range("C9").Value = "Diapositive8"

# The trailing MAXHYP/MINHYP feedback slides were bumped by two slots.
$ws.Range("C15").Value = "Diapositive12"
$ws.Range("C16").Value = "Diapositive14"
$ws.Range("C17").Value = "Diapositive15"
$ws.Range("C18").Value = "Diapositive13"

# Leave the selection where the author left it after the edit.
[void]$ws.Range("C19").Select()
